$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings
# (e.g. "5.938") are not coerced into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.924.20"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.879.47"
$ws.Range("E4").Value = "  +1.45%  "
$ws.Range("D5").Value = "335.33"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").Value = "0.3908"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "46.73"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "0.07924"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "1.006"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "1.914.08"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "5.938"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "7.090"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "0.06794"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "87.40"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "0.00001043"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "16.97"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "1.017"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "27.944.21"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "5.456"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "2.358"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "2.122.46"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "159.44"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "19.90"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").Value = "2.065"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "5.446"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("D31").Value = "120.59"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "0.09519"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "5.313"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "1.346"
$ws.Range("E36").Value = "  -7.36%  "
$ws.Range("D37").Value = "0.06106"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "0.02237"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "1.205"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "8.088"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "0.5866"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "0.1893"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "1.272"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "0.5621"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "12.05"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "3.402"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "0.06855"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "113.35"
$ws.Range("E51").Value = "  +0.77%  "

# Restore default number format / style so unaffected cells
# and the edited cells keep the workbook's original appearance.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
